# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 252
$ws1.Range("F5").Value = 437
$ws1.Range("F7").Value = 567
$ws1.Range("F9").Value = 6838
$ws1.Range("F10").Value = 162
$ws1.Range("F16").Value = 16271
$ws1.Range("F17").Value = 1599
$ws1.Range("F19").Value = 332
$ws1.Range("F21").Value = 117
$ws1.Range("F22").Value = 11403
$ws1.Range("F23").Value = 10
$ws1.Range("F24").Value = 1053
$ws1.Range("F25").Value = 4490
$ws1.Range("F26").Value = 351
$ws1.Range("F29").Value = 846

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 252
$ws4.Range("F5").Value = 437
$ws4.Range("F7").Value = 567
$ws4.Range("F10").Value = 6839
$ws4.Range("F11").Value = 162
$ws4.Range("F18").Value = 16271
$ws4.Range("F19").Value = 1599
$ws4.Range("F21").Value = 332
$ws4.Range("F23").Value = 117
$ws4.Range("F26").Value = 11403
$ws4.Range("F27").Value = 10
$ws4.Range("F28").Value = 1053
$ws4.Range("F29").Value = 4490
$ws4.Range("F30").Value = 351
$ws4.Range("F33").Value = 846
